$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "You are a senior manager in a large retail company, and you are responsible for improving customer experience and increasing online sales. You are tasked with implementing a new digital strategy that includes integrating new technologies into the existing infrastructure and aligning with the overall business goals.What should be your first step?",
        "ques_type": 2,
        "options": [
            "Conduct a thorough review of the current technological landscape to identify any gaps.",
            "Engage with key stakeholders across different departments to gather their insights and input.",
            "Set clear objectives and key performance indicators (KPIs) for the digital strategy.",
            "Evaluate potential risks and challenges associated with implementing new technologies."
        ],
        "score": "Engage with key stakeholders across different departments to gather their insights and input."
    },
    {
        "title": "You are a project manager at a software development company. You are tasked with managing multiple projects simultaneously, with tight deadlines and limited resources. One project involves developing a new mobile app, while another involves enhancing an existing software product. Both projects require immediate attention and are equally important to the company's success.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Request additional resources from senior management to handle both projects simultaneously.",
            "Delegate tasks to your team members based on their skills and expertise.",
            "Assign one project to a senior team member and handle the other one yourself.",
            "Request for an extension of deadlines from both clients to manage the workload."
        ],
        "score": "Delegate tasks to your team members based on their skills and expertise."
    },
    {
        "title": "You are a senior project manager at a large tech company. The company has recently adopted a new project management tool, but you have noticed that some departments are not fully utilizing it. This is impacting the overall efficiency of operations.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Assign a dedicated team member to each department to serve as a point person for questions and support.",
            "Set up regular check-in meetings with each team to review their progress and address any challenges.",
            "Organize additional training sessions on how to use the project management tool.",
            "Create a user-friendly guide or manual that outlines best practices for monitoring progress using the project management tool."
        ],
        "score": "Assign a dedicated team member to each department to serve as a point person for questions and support."
    },
    {
        "title": "You are a project manager in a software development company. Your team is working on a complex project with tight deadlines. One team member is not providing regular updates, making it difficult to track the project's progress.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Schedule a one-on-one meeting with the team member to discuss their progress and any challenges they may be facing.",
            "Request that the team member provide daily progress reports until their tasks are back on track.",
            "Implement a mandatory daily reporting system for all team members.",
            "Increase the frequency of project status meetings to ensure everyone is updated more regularly."
        ],
        "score": "Schedule a one-on-one meeting with the team member to discuss their progress and any challenges they may be facing."
    }
]
'@

$ws.Range("A1:A2").ClearFormats()
$ws.Range("A2").ClearContents()
$ws.Range("A1").Value = $questionsText
$ws.Rows.Item(1).EntireRow.AutoFit()
